# Apply updated dSF (column F) values to reflect repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -3
$ws.Range("F6").Value = -4
$ws.Range("F7").Value = 3
$ws.Range("F9").Value = -6
$ws.Range("F19").Value = 3
$ws.Range("F21").Value = -2
